$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.154.46"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").Value = "2.267.85"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "299.54"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").Value = "95.58"
$ws.Range("E6").Value = "  -3.30%  "

$ws.Range("E7").Value = "  -2.61%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("D10").Value = "33.09"
$ws.Range("E10").Value = "  -4.18%  "

$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").Value = "48.38"
$ws.Range("E12").Value = "  -6.42%  "

$ws.Range("E13").Value = "  +0.75%  "

$ws.Range("D14").Value = "6.65"
$ws.Range("E14").Value = "  -1.53%  "

$ws.Range("D15").Value = "15.60"
$ws.Range("E15").Value = "  -0.67%  "

$ws.Range("D16").Value = "2.622.77"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "2.271.71"
$ws.Range("E17").Value = "  -1.52%  "

$ws.Range("D18").Value = "0.782"
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("D19").Value = "42.144.46"
$ws.Range("E19").Value = "  -0.89%  "

$ws.Range("D20").Value = "11.70"
$ws.Range("E20").Value = "  +1.90%  "

$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").Value = "5.96"
$ws.Range("E22").Value = "  -1.22%  "

$ws.Range("D23").Value = "66.10"
$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("D24").Value = "234.70"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  -2.65%  "

$ws.Range("D28").Value = "23.73"
$ws.Range("E28").Value = "  -5.09%  "

$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -5.19%  "

$ws.Range("D30").Value = "168.10"
$ws.Range("E30").Value = "  +2.60%  "

$ws.Range("D31").Value = "9.16"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("D32").Value = "33.50"
$ws.Range("E32").Value = "  -3.74%  "

$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").Value = "4.86"
$ws.Range("E34").Value = "  -2.91%  "

$ws.Range("D35").Value = "4.56"
$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("D36").Value = "16.62"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("E37").Value = "  -3.34%  "

$ws.Range("D38").Value = "0.0683"
$ws.Range("E38").Value = "  -3.88%  "

$ws.Range("D39").Value = "2.79"
$ws.Range("E39").Value = "  -3.12%  "

$ws.Range("D40").Value = "0.0984"
$ws.Range("E40").Value = "  -1.33%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.108"
$ws.Range("E41").Value = "  -2.72%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "1.71"
$ws.Range("E42").Value = "  -4.91%  "

$ws.Range("E43").Value = "  +1.62%  "

$ws.Range("D44").Value = "1.969.24"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("D45").Value = "0.0276"
$ws.Range("E45").Value = "  -1.25%  "

$ws.Range("D46").Value = "17.39"
$ws.Range("E46").Value = "  -6.27%  "

$ws.Range("D47").Value = "9.52"

$ws.Range("D48").Value = "2.76"
$ws.Range("E48").Value = "  -4.88%  "

$ws.Range("D49").Value = "2.495.15"
$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("D50").Value = "52.18"
$ws.Range("E50").Value = "  -5.71%  "

$ws.Range("D51").Value = "1.47"
$ws.Range("E51").Value = "  -0.95%  "
